$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-19 Saturday" "2025-07-20 Sunday"

Replace-Text "71×53=3763" "29×31=899"
Replace-Text "97×20=1940" "81×99=8019"
Replace-Text "21×42=882" "84×95=7980"
Replace-Text "52×85=4420" "43×33=1419"
Replace-Text "37×64=2368" "42×98=4116"

Replace-Text "51×81=4131" "90×90=8100"
Replace-Text "19×91=1729" "98×74=7252"
Replace-Text "33×51=1683" "96×49=4704"
Replace-Text "29×94=2726" "41×98=4018"
Replace-Text "37×55=2035" "94×83=7802"

Replace-Text "57×23=1311" "34×35=1190"
Replace-Text "13×18=234" "51×43=2193"
Replace-Text "38×37=1406" "76×80=6080"
Replace-Text "46×80=3680" "66×79=5214"
Replace-Text "75×51=3825" "17×57=969"

Replace-Text "65×61=3965" "71×78=5538"
Replace-Text "24×79=1896" "84×52=4368"
Replace-Text "72×34=2448" "19×94=1786"
Replace-Text "61×44=2684" "34×43=1462"
Replace-Text "58×50=2900" "24×19=456"

Replace-Text "26×11=286" "49×76=3724"
Replace-Text "71×84=5964" "31×91=2821"
Replace-Text "64×61=3904" "25×39=975"
Replace-Text "88×33=2904" "80×19=1520"
Replace-Text "63×59=3717" "95×25=2375"
